$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 7710
$ws.Range("J17").Value = 7710
$ws.Range("L17").Value = 23130
$ws.Range("N17").Value = -23466

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 13208238
$ws.Range("I137").Value = 669429.3
$ws.Range("J137").Value = 30306614
$ws.Range("K137").Value = 2008287.9
$ws.Range("L137").Value = 90919842
$ws.Range("M137").Value = -2005737.9
$ws.Range("N137").Value = -90924942

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 20049.477
$ws.Range("I32").Value = 19668.795
$ws.Range("K32").Value = 19668.795
$ws.Range("M32").Value = -19381.795

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 5140.8887
$ws.Range("I61").Value = 5458.5
$ws.Range("K61").Value = 5458.5
$ws.Range("M61").Value = -5246.5

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 17858602
$ws.Range("I74").Value = 19232188
$ws.Range("K74").Value = 19232188
$ws.Range("M74").Value = -19231314

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 17858602
$ws.Range("I77").Value = 19232188
$ws.Range("K77").Value = 96160940
$ws.Range("M77").Value = -96156572

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3516.4092
$ws.Range("I122").Value = 3425.7144
$ws.Range("K122").Value = 10277.1432
$ws.Range("M122").Value = -7827.143199999999

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1866.6279
$ws.Range("I132").Value = 1784.15
$ws.Range("K132").Value = 5352.450000000001
$ws.Range("M132").Value = -2822.450000000001

# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 74343.5
$ws.Range("J135").Value = 74343.5
$ws.Range("L135").Value = 74343.5
$ws.Range("N135").Value = -84483.5

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 5140.8887
$ws.Range("I136").Value = 5458.5
$ws.Range("K136").Value = 16375.5
$ws.Range("M136").Value = -13825.5

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 2238.889
$ws.Range("I99").Value = 2305.5715
$ws.Range("K99").Value = 2305.5715
$ws.Range("M99").Value = -807.5715

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2636.077
$ws.Range("I134").Value = 1929.7084
$ws.Range("J134").Value = 3766.2666
$ws.Range("K134").Value = 5789.1252
$ws.Range("L134").Value = 11298.7998
$ws.Range("M134").Value = -3254.1252
$ws.Range("N134").Value = -16368.7998

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 945.0833
$ws.Range("I16").Value = 829
$ws.Range("K16").Value = 829
$ws.Range("M16").Value = -542

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3123.6365
$ws.Range("I58").Value = 2780.75
$ws.Range("J58").Value = 4038
$ws.Range("K58").Value = 2780.75
$ws.Range("L58").Value = 4038
$ws.Range("M58").Value = -2577.75
$ws.Range("N58").Value = -4444

# Row 113: Patient Patients
$ws.Range("H113").Value = 945.0833
$ws.Range("I113").Value = 829
$ws.Range("K113").Value = 829
$ws.Range("M113").Value = 1341

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 70182340
$ws.Range("I132").Value = 78433110
$ws.Range("K132").Value = 235299330
$ws.Range("M132").Value = -235296800

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2506.2163
$ws.Range("I134").Value = 1667.1154
$ws.Range("K134").Value = 5001.3462
$ws.Range("M134").Value = -2466.3462

# Row 136: Turali Quality
$ws.Range("H136").Value = 3123.6365
$ws.Range("I136").Value = 2780.75
$ws.Range("J136").Value = 4038
$ws.Range("K136").Value = 8342.25
$ws.Range("L136").Value = 12114
$ws.Range("M136").Value = -5792.25
$ws.Range("N136").Value = -17214

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1974.1666
$ws.Range("I97").Value = 1329.8
$ws.Range("K97").Value = 1329.8
$ws.Range("M97").Value = -833.8

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 11115076
$ws.Range("I102").Value = 13517679
$ws.Range("K102").Value = 13517679
$ws.Range("M102").Value = -13516057

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 1275.1428
$ws.Range("I113").Value = 1325.4
$ws.Range("K113").Value = 1325.4
$ws.Range("M113").Value = 844.5999999999999

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 220574.4
$ws.Range("I122").Value = 358545.44
$ws.Range("J122").Value = 5952.8335
$ws.Range("K122").Value = 1075636.32
$ws.Range("L122").Value = 17858.5005
$ws.Range("M122").Value = -1073186.32
$ws.Range("N122").Value = -22758.5005

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2331.258
$ws.Range("I126").Value = 2025.3889
$ws.Range("J126").Value = 2754.7693
$ws.Range("K126").Value = 6076.1667
$ws.Range("L126").Value = 8264.3079
$ws.Range("M126").Value = -3606.1667
$ws.Range("N126").Value = -13204.3079

# Row 132: On Board for Lar
$ws.Range("H132").Value = 107686.266
$ws.Range("I132").Value = 154997.61
$ws.Range("J132").Value = 5178.3335
$ws.Range("K132").Value = 464992.83
$ws.Range("L132").Value = 15535.0005
$ws.Range("M132").Value = -462462.83
$ws.Range("N132").Value = -20595.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 15280522
$ws.Range("I40").Value = 5683606.5
$ws.Range("J40").Value = 41672040
$ws.Range("K40").Value = 5683606.5
$ws.Range("L40").Value = 41672040
$ws.Range("M40").Value = -5683470.5
$ws.Range("N40").Value = -41672312

# Row 57: Too Hot to Handle
$ws.Range("H57").Value = 29166.334
$ws.Range("I57").Value = 24999.5
$ws.Range("J57").Value = 37500
$ws.Range("K57").Value = 24999.5
$ws.Range("L57").Value = 37500
$ws.Range("M57").Value = -24433.5
$ws.Range("N57").Value = -38632

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3034.99
$ws.Range("I132").Value = 2925.8115
$ws.Range("J132").Value = 3278
$ws.Range("K132").Value = 8777.434499999999
$ws.Range("L132").Value = 9834
$ws.Range("M132").Value = -6247.434499999999
$ws.Range("N132").Value = -14894

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4259.5356
$ws.Range("I136").Value = 2827.1292
$ws.Range("K136").Value = 8481.3876
$ws.Range("M136").Value = -5931.3876

$ws = $wb.Worksheets.Item("WVR")
# Row 4: Not Cool Enough
$ws.Range("H4").Value = 833.3333
$ws.Range("I4").Value = 800
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 800
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -687
$ws.Range("N4").Value = -1226

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 10736.5
$ws.Range("I81").Value = 6414.6665
$ws.Range("J81").Value = 17219.25
$ws.Range("K81").Value = 12829.333
$ws.Range("L81").Value = 34438.5
$ws.Range("M81").Value = -11768.333
$ws.Range("N81").Value = -36560.5

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 10736.5
$ws.Range("I84").Value = 6414.6665
$ws.Range("J84").Value = 17219.25
$ws.Range("K84").Value = 64146.665
$ws.Range("L84").Value = 172192.5
$ws.Range("M84").Value = -58842.665
$ws.Range("N84").Value = -182800.5

# Row 96: Skills on Display
$ws.Range("H96").Value = 1896
$ws.Range("J96").Value = 1660.3334
$ws.Range("L96").Value = 1660.3334
$ws.Range("N96").Value = -4406.3334

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 66668932
$ws.Range("I132").Value = 66668932
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 200006796
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -200004266
$ws.Range("N132").ClearContents()

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 3559.5789
$ws.Range("I136").Value = 2526.6897
$ws.Range("K136").Value = 7580.0691
$ws.Range("M136").Value = -5030.0691
